$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.547.17"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "3.021.22"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'509.46"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "'140.45"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.435"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'7.59"
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("D10").Value = "'0.111"
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("D11").Value = "'0.367"
$ws.Range("E11").Value = "  +2.62%  "
$ws.Range("D12").Value = "3.532.81"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "'26.40"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("D16").Value = "57.513.30"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "'6.21"
$ws.Range("E17").Value = "  +3.47%  "
$ws.Range("D18").Value = "3.019.70"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "'12.89"
$ws.Range("E19").Value = "  +2.70%  "
$ws.Range("D20").Value = "'7.99"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").Value = "'328.32"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("D24").Value = "'0.501"
$ws.Range("E24").Value = "  +3.34%  "
$ws.Range("D25").Value = "'64.69"
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "0.0₃0928"
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("D29").Value = "'6.79"
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").Value = "'7.39"
$ws.Range("E30").Value = "  +5.84%  "
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("D32").Value = "'1.20"
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("D33").Value = "'20.67"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("D34").Value = "'4.79"
$ws.Range("E34").Value = "  +4.94%  "
$ws.Range("D35").Value = "'154.35"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "'5.91"
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'24.95"
$ws.Range("E38").Value = "  +4.36%  "
$ws.Range("D39").Value = "'0.0680"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").Value = "3.052.56"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").Value = "'37.94"
$ws.Range("E41").Value = "  +2.57%  "
$ws.Range("D42").Value = "'3.87"
$ws.Range("E42").Value = "  +4.71%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "'0.651"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "2.225.19"
$ws.Range("E46").Value = "  -2.54%  "
$ws.Range("D47").Value = "'0.984"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("E48").Value = "  +3.37%  "
$ws.Range("D49").Value = "'0.0240"
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("D50").Value = "'19.67"
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").Value = "'1.87"
$ws.Range("E51").Value = "  -4.23%  "
